{"js": "// BV_VIP-86: Fixed error in doc UVVM essential mechanisms\n//\n// 1) Body text: \"...shared variables are defined in the UVVM methods package...\"\n//    becomes \"...shared variables are defined in UVVM global signals and shared\n//    variables package...\" (the \"VVC methods package, respectively.\" tail is\n//    untouched).\n// 2) Footer \"Last update\" DATE field cached text: 2018-08-24 -> 2018-11-19.\n\nconst body = context.document.body;\n\nconst textResults = body.search(\n  \"shared variables are defined in the UVVM methods package\",\n  { matchCase: true }\n);\ntextResults.load(\"text\");\nawait context.sync();\n\nif (textResults.items.length > 0) {\n  textResults.items[0].insertText(\n    \"shared variables are defined in UVVM global signals and shared variables package\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const footer = sections.items[i].getFooter(Word.HeaderFooterType.primary);\n  const dateResults = footer.search(\"2018-08-24\", { matchCase: true });\n  dateResults.load(\"text\");\n  await context.sync();\n\n  for (let j = 0; j < dateResults.items.length; j++) {\n    dateResults.items[j].insertText(\"2018-11-19\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# BV_VIP-86: Fixed error in doc UVVM essential mechanisms\n#\n# 1) Body text: \"...shared variables are defined in the UVVM methods package...\"\n#    becomes \"...shared variables are defined in UVVM global signals and shared\n#    variables package...\" (the \"VVC methods package, respectively.\" tail is\n#    untouched).\n# 2) Footer \"Last update\" DATE field cached text: 2018-08-24 -> 2018-11-19.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"shared variables are defined in the UVVM methods package\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"shared variables are defined in UVVM global signals and shared variables package\"\n$find.Execute(\n    [ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n    [ref]$find.Replacement.Text, 2\n)\n\nfor ($i = 1; $i -le $d.Sections.Count; $i++) {\n    $footer = $d.Sections.Item($i).Footers.Item(1)\n    $footerFind = $footer.Range.Find\n    $footerFind.ClearFormatting()\n    $footerFind.Text = \"2018-08-24\"\n    $footerFind.Replacement.ClearFormatting()\n    $footerFind.Replacement.Text = \"2018-11-19\"\n    $footerFind.Execute(\n        [ref]$footerFind.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n        [ref]$footerFind.Replacement.Text, 2\n    )\n}\n"}
